{"js": "// Apply the benchmark-table edits described by the diff.\n// The document contains a single table, one column, many rows; each\n// row's single cell holds one paragraph with (in most rows) a single\n// run of text. A handful of rows (indices 43,44,45 / 1-based 44,45,46)\n// hold a run broken up by <w:tab/> into several <w:t> pieces.\n//\n// Strategy: locate the table, and for each target row, replace the\n// whole paragraph's text via Range.insertText(..., \"Replace\") so the\n// existing run formatting (rFonts/sz) is preserved and any extra runs\n// / tabs in that paragraph collapse into a single run with the new text.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// 0-based row index -> new cell text\nconst edits = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"890\"],\n  [5, \"0.02313\"],\n  [6, \"0.00436\"],\n  [7, \"0.00150\"],\n  [8, \"0.01848\"],\n  [9, \"0.01848\"],\n  [11, \"0.15536\"],\n  [43, \"99.98\"],\n  [44, \"0.16\"],\n  [45, \"766\"],\n];\n\n// Grab first paragraph of each target cell.\nconst paragraphs = [];\nfor (const [rowIndex] of edits) {\n  const cell = table.getCell(rowIndex, 0);\n  const paras = cell.body.paragraphs;\n  paras.load(\"items\");\n  paragraphs.push(paras);\n}\nawait context.sync();\n\n// Replace the text of the first paragraph in each cell.\nedits.forEach(([rowIndex, newText], i) => {\n  const para = paragraphs[i].items[0];\n  const range = para.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-table edits described by the diff.\n# The document has a single, one-column table. Each target row's cell\n# holds one paragraph; most have a single run of text, but a few rows\n# (44, 45, 46) have their run text split across several <w:t> pieces by\n# <w:tab/> characters. Assigning Cell.Range.Text replaces the entire\n# cell content (collapsing any extra runs/tabs) while Word keeps the\n# existing run formatting (rFonts/sz) on the remaining run.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based row number -> new cell text\n$edits = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"890\"\n    6  = \"0.02313\"\n    7  = \"0.00436\"\n    8  = \"0.00150\"\n    9  = \"0.01848\"\n    10 = \"0.01848\"\n    12 = \"0.15536\"\n    44 = \"99.98\"\n    45 = \"0.16\"\n    46 = \"766\"\n}\n\nforeach ($rowNum in $edits.Keys) {\n    $t.Cell($rowNum, 1).Range.Text = $edits[$rowNum]\n}\n"}
